$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Row 7 (EDF Energy (UK)) gets updated data across several columns.
# -----------------------------------------------------------------

# B7 - Annual Revenue (rich text: bold headline, normal detail, bold UK figure, normal tail)
$ws.Range("B7").Value = "€118.7 Billion (Group Sales 2024);`nUK Segment Revenue: ~€13.6 Billion (Est.)."
$ws.Range("B7").Font.Bold = $true
$ws.Range("B7").Characters(15,41).Font.Bold = $false
$ws.Range("B7").Characters(56,14).Font.Bold = $true
$ws.Range("B7").Characters(70,8).Font.Bold = $false

# C7 - Countries covered (rich text: normal lead-in, bold "UK", normal tail)
$ws.Range("C7").Value = "France (HQ), UK, Italy (Edison), Belgium (Luminus);`nActive in 20+ countries worldwide."
$ws.Range("C7").Characters(14,2).Font.Bold = $true
$ws.Range("C7").Characters(16,71).Font.Bold = $false

# D7 - Population of total coverage (rich text: bold headline, normal detail, bold UK figure, normal tail)
$ws.Range("D7").Value = "41.5 Million customer sites globally;`n~3.7 Million in the UK (Residential & Business)."
$ws.Range("D7").Font.Bold = $true
$ws.Range("D7").Characters(14,26).Font.Bold = $false
$ws.Range("D7").Characters(40,11).Font.Bold = $true
$ws.Range("D7").Characters(51,36).Font.Bold = $false

# E7 - Evs (Leasing)
$ws.Range("E7").Value = "Partnership with DriveElectric (Personal & Business leasing); Offers varying makes (Tesla, MG, Nissan, etc.)."

# F7 - EVSE
$ws.Range("F7").Value = "Home: Pod Point Solo 3 (7kW, tethered/untethered).`nPublic: Pod Point Network (Tesco partnership, Lidl, etc.)."

# G7 - BESS
$ws.Range("G7").Value = "Powervault 3 (Eco-Store partnership); Also installs Tesla Powerwall 2 and Sonnen via partners."

# H7 - PV/Solar
$ws.Range("H7").Value = "Installer: EDF Renewables / Hometech.`nHardware: Tier 1 Panels (e.g., Sharp, JA Solar); Inverters: Solis or Growatt."

# I7 - Heat pumps
$ws.Range("I7").Value = "Partnership with CB Heating (EDF acquired them).`nHardware: Primarily Daikin (Altherma 3) and Mitsubishi Electric (Ecodan)."

# Row 7 is now taller to fit the extra detail and B7/D7 become bold-emphasised cells.
$ws.Range("A7:O7").RowHeight = 71.25

# Move the active selection to D7 (Population cell).
$ws.Range("D7").Select() | Out-Null
